$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: set NumberFormat to text before assigning numeric-looking
# strings so Excel stores them as text (matching the source inline strings),
# then reset the style back to Normal so no stray style index is left on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Price / Volume(1h) updates for existing coins ---
Set-TextValue $ws.Range("D2") '25.796.26'
Set-TextValue $ws.Range("E2") '  -0.09%  '
Set-TextValue $ws.Range("D3") '1.638.66'
Set-TextValue $ws.Range("E3") '  +0.34%  '
Set-TextValue $ws.Range("D5") '215.63'
Set-TextValue $ws.Range("E5") '  +0.39%  '
Set-TextValue $ws.Range("E6") '  -0.55%  '
Set-TextValue $ws.Range("E7") '  -0.07%  '
Set-TextValue $ws.Range("E8") '  +0.09%  '
Set-TextValue $ws.Range("D9") '0.0636'
Set-TextValue $ws.Range("E9") '  -0.94%  '
Set-TextValue $ws.Range("D10") '19.63'
Set-TextValue $ws.Range("E10") '  -1.37%  '
Set-TextValue $ws.Range("D11") '0.0793'
Set-TextValue $ws.Range("D12") '4.26'
Set-TextValue $ws.Range("E12") '  +0.35%  '
Set-TextValue $ws.Range("D13") '1.864.18'
Set-TextValue $ws.Range("E13") '  +0.30%  '
Set-TextValue $ws.Range("D14") '1.640.33'
Set-TextValue $ws.Range("E14") '  -1.24%  '
Set-TextValue $ws.Range("D15") '0.564'
Set-TextValue $ws.Range("E15") '  +1.09%  '
Set-TextValue $ws.Range("D17") '63.23'
Set-TextValue $ws.Range("E17") '  +0.30%  '
Set-TextValue $ws.Range("D18") '25.830.74'
Set-TextValue $ws.Range("E18") '  +0.03%  '
Set-TextValue $ws.Range("E19") '  -0.05%  '
Set-TextValue $ws.Range("E20") '  +2.26%  '
Set-TextValue $ws.Range("D21") '192.73'
Set-TextValue $ws.Range("E21") '  -0.54%  '
Set-TextValue $ws.Range("D22") '9.97'
Set-TextValue $ws.Range("D23") '6.29'
Set-TextValue $ws.Range("E23") '  +1.84%  '
Set-TextValue $ws.Range("E24") '  +4.48%  '
Set-TextValue $ws.Range("E25") '  -0.01%  '
Set-TextValue $ws.Range("D26") '141.73'
Set-TextValue $ws.Range("E26") '  +1.28%  '
Set-TextValue $ws.Range("E27") '  +1.23%  '
Set-TextValue $ws.Range("D28") '6.93'
Set-TextValue $ws.Range("E29") '  +0.06%  '
Set-TextValue $ws.Range("E30") '  +0.32%  '
Set-TextValue $ws.Range("E31") '  -0.32%  '
Set-TextValue $ws.Range("E32") '  +0.62%  '
Set-TextValue $ws.Range("E33") '  -0.61%  '
Set-TextValue $ws.Range("E34") '  -0.16%  '
Set-TextValue $ws.Range("E35") '  +0.00%  '
Set-TextValue $ws.Range("D36") '0.906'
Set-TextValue $ws.Range("E36") '  +0.42%  '
Set-TextValue $ws.Range("D37") '1.137.57'
Set-TextValue $ws.Range("E37") '  +1.47%  '
Set-TextValue $ws.Range("E38") '  -1.47%  '
Set-TextValue $ws.Range("D39") '0.545'
Set-TextValue $ws.Range("E39") '  -0.90%  '
Set-TextValue $ws.Range("E41") '  +0.17%  '
Set-TextValue $ws.Range("D42") '5.57'
Set-TextValue $ws.Range("E42") '  +0.76%  '
Set-TextValue $ws.Range("D43") '100.72'
Set-TextValue $ws.Range("E43") '  +1.03%  '
Set-TextValue $ws.Range("D44") '0.803'
Set-TextValue $ws.Range("E44") '  +0.40%  '
Set-TextValue $ws.Range("D45") '1.773.27'
Set-TextValue $ws.Range("E45") '  -0.14%  '

# --- Rows 46-51: list shifted up by one (BabyDogeCoin dropped, Algorand added) ---
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D46") '55.35'
Set-TextValue $ws.Range("E46") '  -0.16%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D47") '0.417'
Set-TextValue $ws.Range("E47") '  -1.15%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D48") '0.0502'
Set-TextValue $ws.Range("E48") '  -0.09%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D49") '1.42'
Set-TextValue $ws.Range("E49") '  +3.61%  '
$ws.Range("B50").Value = 'SynthetixNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextValue $ws.Range("D50") '2.31'
Set-TextValue $ws.Range("E50") '  -0.84%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D51") '0.0956'
Set-TextValue $ws.Range("E51") '  +1.71%  '
